$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.975.64"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "3.387.20"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'571.79"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "'142.13"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.475"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("D11").Value = "'0.387"
$ws.Range("E11").Value = "  -1.13%  "
$ws.Range("D12").Value = "3.964.84"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").Value = "'27.93"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").Value = "3.384.77"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "61.077.83"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("E18").Value = "  -2.82%  "
$ws.Range("E19").Value = "  -3.54%  "
$ws.Range("D20").Value = "'8.95"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("D21").Value = "'384.98"
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("D22").Value = "'75.08"
$ws.Range("E22").Value = "  +2.59%  "
$ws.Range("E23").Value = "  -1.21%  "
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("E25").Value = "  -2.56%  "
$ws.Range("D26").Value = "3.522.67"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "'7.28"
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "'7.97"
$ws.Range("E31").Value = "  -2.00%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'1.38"
$ws.Range("E33").Value = "  -3.90%  "
$ws.Range("D34").Value = "'23.26"
$ws.Range("E34").Value = "  -2.28%  "
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").Value = "'167.13"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").Value = "3.418.49"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "'4.99"
$ws.Range("E38").Value = "  -1.20%  "
$ws.Range("E39").Value = "  -3.40%  "
$ws.Range("D40").Value = "'0.0769"
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("D41").Value = "'26.88"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "'0.780"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").Value = "'4.38"
$ws.Range("E44").Value = "  -1.86%  "
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "2.453.73"
$ws.Range("E47").Value = "  -3.11%  "
$ws.Range("D48").Value = "'22.96"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").Value = "'6.72"
$ws.Range("E49").Value = "  -1.92%  "
$ws.Range("E50").Value = "  +10.37%  "
$ws.Range("E51").Value = "  +1.75%  "
